$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1839.0834
$ws.Range("I135").Value = 1341.5555
$ws.Range("K135").Value = 12073.9995
$ws.Range("M135").Value = -9538.9995

$ws.Range("H138").Value = 3835.319
$ws.Range("I138").Value = 1906.174
$ws.Range("J138").Value = 5684.0835
$ws.Range("K138").Value = 5718.522
$ws.Range("L138").Value = 17052.2505
$ws.Range("M138").Value = -578.5219999999999
$ws.Range("N138").Value = -27332.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3988.535
$ws.Range("J32").Value = 7900
$ws.Range("L32").Value = 7900
$ws.Range("N32").Value = -8474

$ws.Range("H45").Value = 11464.909
$ws.Range("I45").Value = 21250.2
$ws.Range("K45").Value = 21250.2
$ws.Range("M45").Value = -20873.2

$ws.Range("H61").Value = 4458
$ws.Range("I61").Value = 3495
$ws.Range("K61").Value = 3495
$ws.Range("M61").Value = -3283

$ws.Range("H74").Value = 6364
$ws.Range("I74").Value = 1399.5
$ws.Range("J74").Value = 13456.143
$ws.Range("K74").Value = 1399.5
$ws.Range("L74").Value = 13456.143
$ws.Range("M74").Value = -525.5
$ws.Range("N74").Value = -15204.143

$ws.Range("H77").Value = 6364
$ws.Range("I77").Value = 1399.5
$ws.Range("J77").Value = 13456.143
$ws.Range("K77").Value = 6997.5
$ws.Range("L77").Value = 67280.715
$ws.Range("M77").Value = -2629.5
$ws.Range("N77").Value = -76016.715

$ws.Range("H110").Value = 13684.538
$ws.Range("I110").Value = 50999.5
$ws.Range("K110").Value = 50999.5
$ws.Range("M110").Value = -48954.5

$ws.Range("H136").Value = 4458
$ws.Range("I136").Value = 3495
$ws.Range("K136").Value = 10485
$ws.Range("M136").Value = -7935

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 406521.06
$ws.Range("I22").Value = 683.6667
$ws.Range("K22").Value = 683.6667
$ws.Range("M22").Value = -510.6667

$ws.Range("H134").Value = 2545.5757
$ws.Range("I134").Value = 1737.6666
$ws.Range("K134").Value = 5212.9998
$ws.Range("M134").Value = -2677.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 294
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 291
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 291
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = -991

$ws.Range("H31").Value = 63808.766
$ws.Range("I31").Value = 78714.234
$ws.Range("K31").Value = 78714.234
$ws.Range("M31").Value = -78419.234

$ws.Range("H34").Value = 63808.766
$ws.Range("I34").Value = 78714.234
$ws.Range("K34").Value = 78714.234
$ws.Range("M34").Value = -78512.234

$ws.Range("H58").Value = 4046
$ws.Range("I58").Value = 3333.1667
$ws.Range("J58").Value = 6184.5
$ws.Range("K58").Value = 3333.1667
$ws.Range("L58").Value = 6184.5
$ws.Range("M58").Value = -3130.1667
$ws.Range("N58").Value = -6590.5

$ws.Range("H62").Value = 6505.6665
$ws.Range("I62").Value = 5443.875
$ws.Range("K62").Value = 5443.875
$ws.Range("M62").Value = -4819.875

$ws.Range("H64").Value = 47500.332
$ws.Range("J64").Value = 47500.332
$ws.Range("L64").Value = 47500.332
$ws.Range("N64").Value = -47996.332

$ws.Range("H65").Value = 6505.6665
$ws.Range("I65").Value = 5443.875
$ws.Range("K65").Value = 27219.375
$ws.Range("M65").Value = -24099.375

$ws.Range("H67").Value = 47500.332
$ws.Range("J67").Value = 47500.332
$ws.Range("L67").Value = 47500.332
$ws.Range("N67").Value = -49216.332

$ws.Range("H106").Value = 27499.5
$ws.Range("I106").Value = 19999
$ws.Range("J106").Value = 35000
$ws.Range("K106").Value = 19999
$ws.Range("L106").Value = 35000
$ws.Range("M106").Value = -18737
$ws.Range("N106").Value = -37524

$ws.Range("H136").Value = 4046
$ws.Range("I136").Value = 3333.1667
$ws.Range("J136").Value = 6184.5
$ws.Range("K136").Value = 9999.500100000001
$ws.Range("L136").Value = 18553.5
$ws.Range("M136").Value = -7449.500100000001
$ws.Range("N136").Value = -23653.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5032.5
$ws.Range("I5").Value = 509
$ws.Range("J5").Value = 9556
$ws.Range("K5").Value = 1527
$ws.Range("L5").Value = 28668
$ws.Range("M5").Value = -1415
$ws.Range("N5").Value = -28892

$ws.Range("H134").Value = 4669.952
$ws.Range("I134").Value = 2879.375
$ws.Range("K134").Value = 8638.125
$ws.Range("M134").Value = -3568.125

$ws.Range("H135").Value = 5032.5
$ws.Range("I135").Value = 509
$ws.Range("J135").Value = 9556
$ws.Range("K135").Value = 4581
$ws.Range("L135").Value = 86004
$ws.Range("M135").Value = -2046
$ws.Range("N135").Value = -91074

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 222.64706
$ws.Range("I2").Value = 205.04347
$ws.Range("J2").Value = 259.45456
$ws.Range("K2").Value = 205.04347
$ws.Range("L2").Value = 259.45456
$ws.Range("M2").Value = -92.04347000000001
$ws.Range("N2").Value = -485.45456

$ws.Range("H3").Value = 843
$ws.Range("I3").Value = 575.3333
$ws.Range("J3").Value = 1110.6666
$ws.Range("K3").Value = 575.3333
$ws.Range("L3").Value = 1110.6666
$ws.Range("M3").Value = -459.3333
$ws.Range("N3").Value = -1342.6666

$ws.Range("H7").Value = 100000
$ws.Range("J7").Value = 100000
$ws.Range("L7").Value = 100000
$ws.Range("N7").Value = -100224

$ws.Range("H8").Value = 100000
$ws.Range("J8").Value = 100000
$ws.Range("L8").Value = 100000
$ws.Range("N8").Value = -100278

$ws.Range("H9").Value = 374.5
$ws.Range("I9").Value = 300
$ws.Range("J9").Value = 449
$ws.Range("K9").Value = 300
$ws.Range("L9").Value = 449
$ws.Range("M9").Value = -130
$ws.Range("N9").Value = -789

$ws.Range("H10").Value = 667665
$ws.Range("J10").Value = 1497.5
$ws.Range("L10").Value = 1497.5
$ws.Range("N10").Value = -1835.5

$ws.Range("H80").Value = 3996.44
$ws.Range("I80").Value = 2929.4167
$ws.Range("J80").Value = 4981.385
$ws.Range("K80").Value = 2929.4167
$ws.Range("L80").Value = 4981.385
$ws.Range("M80").Value = -1931.4167
$ws.Range("N80").Value = -6977.385

$ws.Range("H83").Value = 3996.44
$ws.Range("I83").Value = 2929.4167
$ws.Range("J83").Value = 4981.385
$ws.Range("K83").Value = 14647.0835
$ws.Range("L83").Value = 24906.925
$ws.Range("M83").Value = -9655.083500000001
$ws.Range("N83").Value = -34890.925

$ws.Range("H132").Value = 628723.5600000001
$ws.Range("I132").Value = 772929.0600000001
$ws.Range("K132").Value = 2318787.18
$ws.Range("M132").Value = -2316257.18

$ws.Range("H135").Value = 50250
$ws.Range("J135").Value = 50250
$ws.Range("L135").Value = 50250
$ws.Range("N135").Value = -60390

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6989.88
$ws.Range("I7").Value = 8124.8335
$ws.Range("K7").Value = 8124.8335
$ws.Range("M7").Value = -8012.8335

$ws.Range("H55").Value = 422.5
$ws.Range("I55").Value = 554
$ws.Range("J55").Value = 159.5
$ws.Range("K55").Value = 554
$ws.Range("L55").Value = 159.5
$ws.Range("M55").Value = -381
$ws.Range("N55").Value = -505.5

$ws.Range("H68").Value = 20432.727
$ws.Range("I68").Value = 4090
$ws.Range("J68").Value = 29771.428
$ws.Range("K68").Value = 4090
$ws.Range("L68").Value = 29771.428
$ws.Range("M68").Value = -3341
$ws.Range("N68").Value = -31269.428

$ws.Range("H71").Value = 20432.727
$ws.Range("I71").Value = 4090
$ws.Range("J71").Value = 29771.428
$ws.Range("K71").Value = 20450
$ws.Range("L71").Value = 148857.14
$ws.Range("M71").Value = -16706
$ws.Range("N71").Value = -156345.14

$ws.Range("H126").Value = 6989.88
$ws.Range("I126").Value = 8124.8335
$ws.Range("K126").Value = 24374.5005
$ws.Range("M126").Value = -21904.5005

$ws.Range("H132").Value = 4749.353
$ws.Range("I132").Value = 3477.9167
$ws.Range("K132").Value = 10433.7501
$ws.Range("M132").Value = -7903.750100000001

$ws.Range("H136").Value = 4681.121
$ws.Range("I136").Value = 4541.136
$ws.Range("J136").Value = 4961.091
$ws.Range("K136").Value = 13623.408
$ws.Range("L136").Value = 14883.273
$ws.Range("M136").Value = -11073.408
$ws.Range("N136").Value = -19983.273

$ws.Range("H139").Value = 98000
$ws.Range("J139").Value = 98000
$ws.Range("L139").Value = 98000
$ws.Range("N139").Value = -108280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 49999.125
$ws.Range("I2").Value = 49999.75
$ws.Range("J2").Value = 49998.5
$ws.Range("K2").Value = 49999.75
$ws.Range("L2").Value = 49998.5
$ws.Range("M2").Value = -49887.75
$ws.Range("N2").Value = -50222.5

$ws.Range("H136").Value = 3240.5
$ws.Range("I136").Value = 3302.923
$ws.Range("J136").Value = 2970
$ws.Range("K136").Value = 9908.769
$ws.Range("L136").Value = 8910
$ws.Range("M136").Value = -7358.769
$ws.Range("N136").Value = -14010
